$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")

# --- Block 1 (rows 19-25): fill in the "W20" week column (Y) ---
$ws.Range("Y19").Value = 4
$ws.Range("Y20").Value = 0

# --- Block 4 (rows 67-76): fill in / correct the "Y" week column ---
$ws.Range("Y67").Value = 3
$ws.Range("Y68").Value = 10
$ws.Range("Y69").Value = 8
$ws.Range("Y70").Value = 3
$ws.Range("X71").Value = 8
$ws.Range("Y71").Value = 8

# --- Block 5 (rows 85-94): fill in the last two week columns ---
$ws.Range("Y86").Value = 5
$ws.Range("Y87").Value = 5
$ws.Range("X88").Value = 1
$ws.Range("Y88").Value = 1
$ws.Range("Y89").Value = 1
$ws.Range("V90").Value = 1
$ws.Range("W90").Value = 1

# --- New grand-total rows under the last block ---
$ws.Range("Y115").Value = "Total Plan"
$ws.Range("AA115").Formula = "=SUM(AA14,AA28,AA45,AA62,AA80,AA98)"

$ws.Range("Y116").Value = "Total Real"
$ws.Range("AA116").Formula = "=SUM(AA22,AA56,AA39,AA73,AA91,AA109)"

# --- View state: move the selection on "Kosten" ... ---
$ws.Range("Y117").Select()

# ... and make "Status-4" the active/selected sheet, as in the target file.
$ws4 = $wb.Worksheets.Item("Status-4")
$ws4.Activate()
